$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bulk single-cell numeric corrections across rows 98-639
$ws.Range("F98").Value = 13886949
$ws.Range("G98").Value = 11076382
$ws.Range("M98").Value = 8643090
$ws.Range("M104").Value = 8353740
$ws.Range("M108").Value = 8263904
$ws.Range("D140").Value = 823714
$ws.Range("N140").Value = 8935550
$ws.Range("N141").Value = 8782065
$ws.Range("N142").Value = 8776333
$ws.Range("N143").Value = 8836775
$ws.Range("N144").Value = 9042583
$ws.Range("N145").Value = 9395203
$ws.Range("N146").Value = 9716067
$ws.Range("N147").Value = 9916898
$ws.Range("N148").Value = 10176531
$ws.Range("N149").Value = 10402112
$ws.Range("N150").Value = 10586512
$ws.Range("N151").Value = 10834425
$ws.Range("J156").Value = 37.93
$ws.Range("D191").Value = 573564
$ws.Range("E191").Value = 207671
$ws.Range("G191").Value = 8670689
$ws.Range("H191").Value = 1073112
$ws.Range("M191").Value = 14527950
$ws.Range("N191").Value = 8368524
$ws.Range("O191").Value = 2551082
$ws.Range("M192").Value = 14532725
$ws.Range("N192").Value = 8290098
$ws.Range("O192").Value = 2547704
$ws.Range("M193").Value = 14761249
$ws.Range("N193").Value = 8365190
$ws.Range("O193").Value = 2572352
$ws.Range("F194").Value = 14729044
$ws.Range("M194").Value = 14843930
$ws.Range("N194").Value = 8388077
$ws.Range("O194").Value = 2570152
$ws.Range("M195").Value = 14888860
$ws.Range("N195").Value = 8325549
$ws.Range("O195").Value = 2571041
$ws.Range("F196").Value = 13097747
$ws.Range("M196").Value = 14967542
$ws.Range("N196").Value = 8320045
$ws.Range("O196").Value = 2563115
$ws.Range("M197").Value = 15089686
$ws.Range("N197").Value = 8341509
$ws.Range("O197").Value = 2566570
$ws.Range("M198").Value = 15141561
$ws.Range("N198").Value = 8330324
$ws.Range("O198").Value = 2567371
$ws.Range("F199").Value = 10187293
$ws.Range("M199").Value = 15309105
$ws.Range("N199").Value = 8356572
$ws.Range("O199").Value = 2587921
$ws.Range("M200").Value = 15440714
$ws.Range("N200").Value = 8350591
$ws.Range("O200").Value = 2599972
$ws.Range("M201").Value = 15510489
$ws.Range("N201").Value = 8368334
$ws.Range("O201").Value = 2605545
$ws.Range("F202").Value = 10071161
$ws.Range("L202").Value = 134.04
$ws.Range("M202").Value = 15673799
$ws.Range("N202").Value = 8381891
$ws.Range("O202").Value = 2635462
$ws.Range("F203").Value = 9851303
$ws.Range("B232").Value = 1361125
$ws.Range("F232").Value = 11588933
$ws.Range("D622").Value = 478045
$ws.Range("E622").Value = 305612
$ws.Range("H622").Value = 2683520
$ws.Range("K622").Value = 20.98
$ws.Range("L622").Value = 340.1
$ws.Range("M622").Value = 58985486
$ws.Range("N622").Value = 8911033
$ws.Range("O622").Value = 8847971
$ws.Range("J623").Value = 34.87
$ws.Range("K623").Value = 19.46
$ws.Range("L623").Value = 344.45
$ws.Range("M623").Value = 52913702
$ws.Range("N623").Value = 8503869
$ws.Range("O623").Value = 6719966
$ws.Range("J624").Value = 39.27
$ws.Range("L624").Value = 347.33
$ws.Range("M624").Value = 48456932
$ws.Range("N624").Value = 8091688
$ws.Range("O624").Value = 5601360
$ws.Range("J625").Value = 40.31
$ws.Range("K625").Value = 17.8
$ws.Range("L625").Value = 349.78
$ws.Range("M625").Value = 44259956
$ws.Range("N625").Value = 7631462
$ws.Range("O625").Value = 4741942
$ws.Range("J626").Value = 40.18
$ws.Range("L626").Value = 352.67
$ws.Range("M626").Value = 41276761
$ws.Range("N626").Value = 7250472
$ws.Range("O626").Value = 4097878
$ws.Range("J627").Value = 38.57
$ws.Range("L627").Value = 355.82
$ws.Range("M627").Value = 38624133
$ws.Range("N627").Value = 6876981
$ws.Range("O627").Value = 3603849
$ws.Range("J628").Value = 35.5
$ws.Range("L628").Value = 359.78
$ws.Range("M628").Value = 35818386
$ws.Range("N628").Value = 6397086
$ws.Range("O628").Value = 3163538
$ws.Range("J629").Value = 33.87
$ws.Range("L629").Value = 363.43
$ws.Range("M629").Value = 33636671
$ws.Range("N629").Value = 5635512
$ws.Range("O629").Value = 2880185
$ws.Range("J630").Value = 32.57
$ws.Range("L630").Value = 367.89
$ws.Range("M630").Value = 31534221
$ws.Range("N630").Value = 5196379
$ws.Range("O630").Value = 2635215
$ws.Range("C631").Value = 1297.2
$ws.Range("J631").Value = 31.55
$ws.Range("L631").Value = 372.77
$ws.Range("M631").Value = 29240156
$ws.Range("N631").Value = 4865258
$ws.Range("O631").Value = 2407594
$ws.Range("J632").Value = 30.82
$ws.Range("L632").Value = 377.44
$ws.Range("M632").Value = 27537822
$ws.Range("N632").Value = 4641621
$ws.Range("O632").Value = 2234803
$ws.Range("J633").Value = 30.24
$ws.Range("L633").Value = 382.69
$ws.Range("M633").Value = 26108626
$ws.Range("N633").Value = 4537767
$ws.Range("O633").Value = 2079434
$ws.Range("J634").Value = 29.34
$ws.Range("J635").Value = 29.9
$ws.Range("J636").Value = 30.52
$ws.Range("J637").Value = 31.21
$ws.Range("J638").Value = 31.91
$ws.Range("J639").Value = 32.19
$ws.Range("R649").Value = 446894
$ws.Range("P650").Value = 5420
$ws.Range("Q650").Value = 14742
$ws.Range("R650").Value = 443767

# Row 651 corrections + new cells C651, R651
# Row 652 is a brand-new month appended at the bottom
$ws.Range("B651").Value = 881571
$ws.Range("C651").Value = 2059.2
$ws.Range("D651").Value = 487456
$ws.Range("E651").Value = 157472
$ws.Range("F651").Value = 8524590
$ws.Range("G651").Value = 7385863
$ws.Range("H651").Value = 3339684
$ws.Range("I651").Value = 459.48
$ws.Range("J651").Value = 37.21
$ws.Range("K651").Value = 14.99
$ws.Range("L651").Value = 443.77
$ws.Range("M651").Value = 33854725
$ws.Range("N651").Value = 5182038
$ws.Range("O651").Value = 1855331
$ws.Range("P651").Value = 5431
$ws.Range("Q651").Value = 14784
$ws.Range("R651").Value = 438817
$ws.Range("A652").Value = 45382
$ws.Range("B652").Value = 1180153
$ws.Range("D652").Value = 451949
$ws.Range("E652").Value = 144376
$ws.Range("F652").Value = 8854339
$ws.Range("G652").Value = 7495292
$ws.Range("H652").Value = 3308275
$ws.Range("I652").Value = 448.03
$ws.Range("J652").Value = 37.16
$ws.Range("K652").Value = 14.99
$ws.Range("L652").Value = 444.41
$ws.Range("M652").Value = 34186103
$ws.Range("N652").Value = 5223413
$ws.Range("O652").Value = 1869259
$ws.Range("P652").Value = 5444
$ws.Range("Q652").Value = 14833
